# "Migliorie interfaccia e messaggi"
#
# 1) UI/window tweak: resize the workbook window (reflected, when the host
#    supports it, as xWindow/yWindow/windowWidth/windowHeight on the saved
#    <workbookView>).
# 2) Message tweak: rename two placeholder test strings ("Pippo" -> "Caso 1",
#    "Pluto" -> "Caso D") used as lookup-table keys in columns Y and Z.
# 3) The user then re-sorted column Y (ascending) and column Z (descending,
#    with header) independently, leaving Excel's remembered sort state
#    behind, and left the Y4:Z11 block selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) window geometry -----------------------------------------------
try {
    $win = $wb.Windows.Item(1)
    $win.Left   = 2295
    $win.Top    = 2295
    $win.Width  = 28800
    $win.Height = 15345
} catch {
    # Window sizing isn't modeled everywhere; ignore if unsupported.
}

# --- 2) rename the two placeholder labels, in place --------------------
# "Pluto" (Y11/Z11) becomes "Caso D" ...
$ws.Range("Y11").Value = "Caso D"
$ws.Range("Z11").Value = "Caso D"
# ... and "Pippo" (Y7/Z7) becomes "Caso 1"
$ws.Range("Y7").Value = "Caso 1"
$ws.Range("Z7").Value = "Caso 1"

# --- 3a) sort column Y (Y4:Y11) ascending, on its own -------------------
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("Y4:Y11"), 0, 1, 0, 0)
$ws.Sort.SetRange($ws.Range("Y4:Y11"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# --- 3b) sort column Z (Z4:Z11) descending, with header ------------------
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("Z4:Z11"), 0, 2, 0, 0)
$ws.Sort.SetRange($ws.Range("Z4:Z11"))
$ws.Sort.Header = 1
$ws.Sort.Apply()

# --- 3c) leave the block selected, as the author did ---------------------
$ws.Range("Y4:Z11").Select() | Out-Null
